$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 3785.7144
$ws.Range("J7").Value = 3785.7144
$ws.Range("L7").Value = 3785.7144
$ws.Range("N7").Value = -4009.7144

$ws.Range("H14").Value = 3785.7144
$ws.Range("J14").Value = 3785.7144
$ws.Range("L14").Value = 3785.7144
$ws.Range("N14").Value = -4167.7144

$ws.Range("H92").Value = 7692670
$ws.Range("I92").Value = 8333684
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 8333684
$ws.Range("L92").Value = 500
$ws.Range("M92").Value = -8332436
$ws.Range("N92").Value = -2996

$ws.Range("H106").Value = 4622.25
$ws.Range("I106").Value = 4622.25
$ws.Range("K106").Value = 4622.25
$ws.Range("M106").Value = -3991.25

$ws.Range("H138").Value = 1802.4517
$ws.Range("J138").Value = 2496.7273
$ws.Range("L138").Value = 7490.1819
$ws.Range("N138").Value = -17770.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5951.2407
$ws.Range("I32").Value = 2266.7046
$ws.Range("J32").Value = 22163.2
$ws.Range("K32").Value = 2266.7046
$ws.Range("L32").Value = 22163.2
$ws.Range("M32").Value = -1979.7046
$ws.Range("N32").Value = -22737.2

$ws.Range("H61").Value = 79753.46000000001
$ws.Range("I61").Value = 2755.875
$ws.Range("J61").Value = 202949.6
$ws.Range("K61").Value = 2755.875
$ws.Range("L61").Value = 202949.6
$ws.Range("M61").Value = -2543.875
$ws.Range("N61").Value = -203373.6

$ws.Range("H74").Value = 94269.55
$ws.Range("I74").Value = 167783.5
$ws.Range("K74").Value = 167783.5
$ws.Range("M74").Value = -166909.5

$ws.Range("H77").Value = 94269.55
$ws.Range("I77").Value = 167783.5
$ws.Range("K77").Value = 838917.5
$ws.Range("M77").Value = -834549.5

$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490

$ws.Range("H122").Value = 4128.9165
$ws.Range("I122").Value = 2756.125
$ws.Range("J122").Value = 6874.5
$ws.Range("K122").Value = 8268.375
$ws.Range("L122").Value = 20623.5
$ws.Range("M122").Value = -5818.375
$ws.Range("N122").Value = -25523.5

$ws.Range("H136").Value = 79753.46000000001
$ws.Range("I136").Value = 2755.875
$ws.Range("J136").Value = 202949.6
$ws.Range("K136").Value = 8267.625
$ws.Range("L136").Value = 608848.8
$ws.Range("M136").Value = -5717.625
$ws.Range("N136").Value = -613948.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1295.2667
$ws.Range("J20").Value = 1373.25
$ws.Range("L20").Value = 1373.25
$ws.Range("N20").Value = -1867.25

$ws.Range("H42").Value = 350000
$ws.Range("J42").Value = 350000
$ws.Range("L42").Value = 350000
$ws.Range("N42").Value = -350656

$ws.Range("H86").Value = 4326.421
$ws.Range("I86").Value = 3904.9
$ws.Range("K86").Value = 3904.9
$ws.Range("M86").Value = -2781.9

$ws.Range("H89").Value = 4326.421
$ws.Range("I89").Value = 3904.9
$ws.Range("K89").Value = 19524.5
$ws.Range("M89").Value = -13908.5

$ws.Range("H132").Value = 34451.906
$ws.Range("J132").Value = 34451.906
$ws.Range("L132").Value = 34451.906
$ws.Range("N132").Value = -44571.906

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1642.409
$ws.Range("I58").Value = 1487.0588
$ws.Range("K58").Value = 1487.0588
$ws.Range("M58").Value = -1284.0588

$ws.Range("H132").Value = 1517504.4
$ws.Range("I132").Value = 1685754.9
$ws.Range("K132").Value = 5057264.699999999
$ws.Range("M132").Value = -5054734.699999999

$ws.Range("H134").Value = 2265760
$ws.Range("I134").Value = 2553421.8
$ws.Range("J134").Value = 252127
$ws.Range("K134").Value = 7660265.399999999
$ws.Range("L134").Value = 756381
$ws.Range("M134").Value = -7657730.399999999
$ws.Range("N134").Value = -761451

$ws.Range("H136").Value = 1642.409
$ws.Range("I136").Value = 1487.0588
$ws.Range("K136").Value = 4461.1764
$ws.Range("M136").Value = -1911.1764

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 2536
$ws.Range("J121").Value = 3123.375
$ws.Range("L121").Value = 9370.125
$ws.Range("N121").Value = -11990.125

$ws.Range("H131").Value = 1568.4445
$ws.Range("I131").Value = 1020.55554
$ws.Range("J131").Value = 2116.3333
$ws.Range("K131").Value = 3061.66662
$ws.Range("L131").Value = 6348.999899999999
$ws.Range("M131").Value = 1978.33338
$ws.Range("N131").Value = -16428.9999

$ws.Range("H137").Value = 4119.4707
$ws.Range("I137").Value = 2239.8
$ws.Range("J137").Value = 6804.7144
$ws.Range("K137").Value = 6719.400000000001
$ws.Range("L137").Value = 20414.1432
$ws.Range("M137").Value = -1619.400000000001
$ws.Range("N137").Value = -30614.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 13999.75
$ws.Range("I7").Value = 11500
$ws.Range("J7").Value = 16499.5
$ws.Range("K7").Value = 11500
$ws.Range("L7").Value = 16499.5
$ws.Range("M7").Value = -11388
$ws.Range("N7").Value = -16723.5

$ws.Range("H8").Value = 13999.75
$ws.Range("I8").Value = 11500
$ws.Range("J8").Value = 16499.5
$ws.Range("K8").Value = 11500
$ws.Range("L8").Value = 16499.5
$ws.Range("M8").Value = -11361
$ws.Range("N8").Value = -16777.5

$ws.Range("H12").Value = 5000
$ws.Range("I12").Value = 5000
$ws.Range("K12").Value = 5000
$ws.Range("M12").Value = -4860

$ws.Range("H14").Value = 23753500
$ws.Range("I14").Value = 47500000
$ws.Range("J14").Value = 7000
$ws.Range("K14").Value = 47500000
$ws.Range("L14").Value = 7000
$ws.Range("M14").Value = -47499832
$ws.Range("N14").Value = -7336

$ws.Range("H19").Value = 16999.5
$ws.Range("I19").Value = 15499.25
$ws.Range("J19").Value = 20000
$ws.Range("K19").Value = 15499.25
$ws.Range("L19").Value = 20000
$ws.Range("M19").Value = -15211.25
$ws.Range("N19").Value = -20576

$ws.Range("H59").Value = 2100
$ws.Range("J59").Value = 2950
$ws.Range("L59").Value = 2950
$ws.Range("N59").Value = -4116

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H126").Value = 3401.1
$ws.Range("I126").Value = 2388.6667
$ws.Range("J126").Value = 4229.4546
$ws.Range("K126").Value = 7166.000100000001
$ws.Range("L126").Value = 12688.3638
$ws.Range("M126").Value = -4696.000100000001
$ws.Range("N126").Value = -17628.3638

$ws.Range("H51").Value = 15000
$ws.Range("J51").Value = 15000
$ws.Range("L51").Value = 15000
$ws.Range("N51").Value = -15956

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1001.3333
$ws.Range("I61").Value = 1001.3333
$ws.Range("K61").Value = 1001.3333
$ws.Range("M61").Value = -799.3333

$ws.Range("H100").Value = 4366.533
$ws.Range("I100").Value = 3999.8462
$ws.Range("J100").Value = 6750
$ws.Range("K100").Value = 3999.8462
$ws.Range("L100").Value = 6750
$ws.Range("M100").Value = -3458.8462
$ws.Range("N100").Value = -7832

$ws.Range("H113").Value = 1001.3333
$ws.Range("I113").Value = 1001.3333
$ws.Range("K113").Value = 1001.3333
$ws.Range("M113").Value = 1168.6667

$ws.Range("H132").Value = 3842.5
$ws.Range("I132").Value = 3842.5
$ws.Range("K132").Value = 11527.5
$ws.Range("M132").Value = -8997.5

$ws.Range("H134").Value = 139673.25
$ws.Range("J134").Value = 139673.25
$ws.Range("L134").Value = 139673.25
$ws.Range("N134").Value = -149813.25

$ws.Range("H136").Value = 1862.05
$ws.Range("I136").Value = 1589.5
$ws.Range("K136").Value = 4768.5
$ws.Range("M136").Value = -2218.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 42629.8
$ws.Range("I101").Value = 20547
$ws.Range("J101").Value = 48150.5
$ws.Range("K101").Value = 20547
$ws.Range("L101").Value = 48150.5
$ws.Range("N101").Value = -54640.5
$ws.Range("M101").Value = -17302

$ws.Range("H132").Value = 2830.6924
$ws.Range("I132").Value = 2599.875
$ws.Range("K132").Value = 7799.625
$ws.Range("M132").Value = -5269.625
